$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "Recorded By" (column G) list must NOT be reordered.
$excludedRows = @(7, 33, 59)

$lastRow = 157
for ($row = 2; $row -le $lastRow; $row++) {
    if ($excludedRows -contains $row) {
        continue
    }

    $cell = $ws.Cells.Item($row, 7)   # column G = "Recorded By"
    $value = $cell.Value2

    if ($null -eq $value) {
        continue
    }

    $parts = $value -split ", "
    if ($parts.Count -gt 1) {
        $reversed = $parts[($parts.Count - 1)..0]
        $cell.Value = [string]::Join(", ", $reversed)
    }
}
